# Updates cryptos list prices/volume figures per upstream commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.797.01"
$ws.Range("E2").Value = "  -1.33%  "
$ws.Range("D3").Value = "2.907.91"
$ws.Range("E3").Value = "  -1.83%  "
$ws.Range("D4").Formula = '="1.00"'
$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Formula = '="586.58"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -1.51%  "
$ws.Range("D6").Formula = '="146.39"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "2.907.89"
$ws.Range("E9").Value = "  -1.81%  "
$ws.Range("D10").Formula = '="6.83"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  -7.11%  "
$ws.Range("E11").Value = "  +5.28%  "
$ws.Range("D12").Formula = '="0.434"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  -3.05%  "
$ws.Range("D13").Formula = '="0.0000236"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("D14").Formula = '="32.74"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  -1.88%  "
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("E16").Value = "  -1.69%  "
$ws.Range("D17").Value = "61.851.32"
$ws.Range("E17").Value = "  -1.12%  "
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("D19").Value = "2.909.06"
$ws.Range("E19").Value = "  -2.75%  "
$ws.Range("D20").Formula = '="435.80"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("D21").Formula = '="13.36"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -0.53%  "
$ws.Range("E22").Value = "  -1.87%  "
$ws.Range("D23").Formula = '="6.95"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  -1.96%  "
$ws.Range("D24").Formula = '="80.55"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -1.65%  "
$ws.Range("D25").Formula = '="11.96"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").Formula = '="10.21"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  -9.03%  "
$ws.Range("D27").Formula = '="2.07"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  -2.84%  "
$ws.Range("E29").Value = "  +21.84%  "
$ws.Range("D30").Formula = '="7.14"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("E31").Value = "  -1.89%  "
$ws.Range("E32").Value = "  -0.72%  "
$ws.Range("D33").Formula = '="0.108"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").Formula = '="25.87"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  -2.92%  "
$ws.Range("D36").Formula = '="0.974"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  -1.65%  "
$ws.Range("D37").Formula = '="5.50"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  -2.59%  "
$ws.Range("E38").Value = "  +3.81%  "
$ws.Range("D39").Formula = '="49.12"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  -1.05%  "
$ws.Range("E40").Value = "  -2.28%  "
$ws.Range("D41").Formula = '="8.39"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  -2.02%  "
$ws.Range("E42").Value = "  -1.58%  "
$ws.Range("D43").Formula = '="0.272"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -3.52%  "
$ws.Range("D44").Formula = '="38.87"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("D45").Value = "2.698.34"
$ws.Range("E45").Value = "  -0.77%  "
$ws.Range("D46").Formula = '="133.97"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  -1.47%  "
$ws.Range("D47").Formula = '="0.0335"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -1.86%  "
$ws.Range("D48").Formula = '="341.98"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  -5.90%  "
$ws.Range("E50").Value = "  -1.80%  "
$ws.Range("D51").Formula = '="22.20"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  -3.07%  "
